$d = $word.ActiveDocument

$pairs = @(
    @("53÷6=", "57÷9="),
    @("69÷3=", "94÷4="),
    @("17÷7=", "88÷9="),
    @("45÷4=", "45÷2="),
    @("59÷8=", "59÷2="),
    @("28÷5=", "36÷6="),
    @("17÷6=", "75÷8="),
    @("58÷7=", "35÷9="),
    @("35÷5=", "69÷9="),
    @("42÷7=", "84÷9="),
    @("12÷9=", "70÷7="),
    @("55÷5=", "50÷2="),
    @("23÷4=", "96÷9="),
    @("80÷9=", "57÷9="),
    @("46÷2=", "51÷9="),
    @("21÷3=", "31÷8="),
    @("93÷6=", "18÷7="),
    @("54÷3=", "82÷7="),
    @("49÷9=", "81÷8="),
    @("71÷7=", "56÷8="),
    @("15÷2=", "10÷7="),
    @("60÷8=", "11÷6="),
    @("46÷6=", "70÷6="),
    @("68÷3=", "75÷9="),
    @("20÷6=", "66÷8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
